$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is a grid of election-simulation results:
#   Row 1: "nCandidates:" header, one value per block (merged across block)
#   Row 2: "nVoters:" header, repeating series of values inside each block
#   Rows 4-12: one method per row, winner name repeated across all columns
#
# Old layout: 6 blocks of 8 columns (nCandidates = 3,5,8,10,12,14),
#             nVoters series per block = 3,5,10,50,100,500,1000,10000
# New layout: 7 blocks of 6 columns (nCandidates = 2,4,6,8,10,12,14),
#             nVoters series per block = 10,50,100,500,1000,10000
#
# Strategy: reshape the grid using ONLY whole-column insert/delete so the
# existing merged ranges (and their styles) get resized/shifted by Excel
# automatically, instead of calling Merge()/UnMerge() (which would mint a
# brand-new, spurious cell style for every merge operation).
# ---------------------------------------------------------------------------

# Step 1: drop the two blocks that disappear entirely (old nCandidates 3 and
# 5), i.e. columns B:Q (16 columns). The remaining 4 blocks (8,10,12,14)
# slide left to start at column B.
$ws.Range("B1:Q12").EntireColumn.Delete() | Out-Null

# Step 2: each surviving block shrinks from 8 columns to 6 (the nVoters
# series loses its leading "3" and "5" entries). Work right-to-left so
# earlier deletes don't invalidate the column letters used later.
$ws.Range("AF1:AG12").EntireColumn.Delete() | Out-Null
$ws.Range("X1:Y12").EntireColumn.Delete() | Out-Null
$ws.Range("P1:Q12").EntireColumn.Delete() | Out-Null
$ws.Range("H1:I12").EntireColumn.Delete() | Out-Null

# Step 3: insert 3 fresh 6-column blocks before column B for the new
# nCandidates values 2, 4, 6. Inserting the same width 3 times in a row at
# the same spot pushes everything already there further right each time.
$ws.Range("B1:G12").EntireColumn.Insert() | Out-Null
$ws.Range("B1:G12").EntireColumn.Insert() | Out-Null
$ws.Range("B1:G12").EntireColumn.Insert() | Out-Null

# At this point the sheet spans A1:AQ12 with 7 correctly-sized, still-merged
# header blocks at B,H,N,T,Z,AF,AL - the rightmost four (T,Z,AF,AL) already
# hold the right header/body data (8,10,12,14); B,H,N are blank new columns
# that still need values.

$blockStarts = @("B", "H", "N", "T", "Z", "AF", "AL")
$nCandidates = @(2, 4, 6, 8, 10, 12, 14)
$nVoters     = @(10, 50, 100, 500, 1000, 10000)

for ($i = 0; $i -lt $blockStarts.Count; $i++) {
    $startIdx = $ws.Range($blockStarts[$i] + "1").Column

    # Row 1: nCandidates value lives only in the first column of the block
    $ws.Cells.Item(1, $startIdx).Value = $nCandidates[$i]

    # Row 2: nVoters series repeats across the 6 columns of the block
    for ($j = 0; $j -lt $nVoters.Count; $j++) {
        $ws.Cells.Item(2, $startIdx + $j).Value = $nVoters[$j]
    }

    # Rows 4-12: default winner is "Montebourg" in every data column
    for ($r = 4; $r -le 12; $r++) {
        for ($j = 0; $j -lt 6; $j++) {
            $ws.Cells.Item($r, $startIdx + $j).Value = "Montebourg"
        }
    }
}

# ---------------------------------------------------------------------------
# Specific winner corrections (bug fix highlighted by the commit message):
# some configurations actually elect "Jadot" or "Philipot" instead of the
# default "Montebourg".
# ---------------------------------------------------------------------------
foreach ($ref in @("Z5", "AF5", "AL5", "Z6", "AF6", "AL6")) {
    $ws.Range($ref).Value = "Jadot"
}

foreach ($ref in @("U7", "T11", "Z11", "AF11", "AL11")) {
    $ws.Range($ref).Value = "Philipot"
}
